$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44; this shifts the existing rows 44-149
# down to 45-150 (carrying their values/formatting with them), and the
# sheet's used range grows to A1:R150 automatically.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new data record.
$ws.Cells.Item(44, 1).Value = 6
$ws.Cells.Item(44, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(44, 3).Value = "Metropolitana"

# Column D carries the same custom date/time number format used by the
# rest of the column (style index 2 in the original workbook).
$ws.Cells.Item(44, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(44, 4).Value = 44560

$ws.Cells.Item(44, 5).Value = 13
$ws.Cells.Item(44, 6).Value = 100112001
$ws.Cells.Item(44, 7).Value = "Berenjena"
$ws.Cells.Item(44, 8).Value = "Sin especificar"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 400
$ws.Cells.Item(44, 11).Value = 8000
$ws.Cells.Item(44, 12).Value = 9000
$ws.Cells.Item(44, 13).Value = 8575
$ws.Cells.Item(44, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(44, 15).Value = "Región Metropolitana"
$ws.Cells.Item(44, 16).Value = 143
$ws.Cells.Item(44, 17).Value = 60
$ws.Cells.Item(44, 18).Value = "Hortaliza"
